$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$oldText = "'Hispanic or Latino;White', 'Hispanic or Latino'"
$newText = "'Hispanic or Latino;White', 'Hispanic or Latino', 'Black or African American;Hispanic or Latino', 'Asian;Hispanic or Latino'"

$cells = @("B2", "C2", "B3", "B4", "B5", "B6", "B7")

foreach ($addr in $cells) {
    $rng = $ws.Range($addr)
    $text = $rng.Value2
    $updated = $text.Replace($oldText, $newText)
    $rng.Value = $updated
}

Write-Host "Done updating cells"
